$d = $word.ActiveDocument

function Replace-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    # Trim the trailing paragraph mark from the range so Find/Replace only
    # touches the visible text of the paragraph.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

# --- 1:1 text replacements for the first 20 paragraphs ---------------------
Replace-ParaText 1  "OCR Results - diplome licence allemand.pdf"
Replace-ParaText 2  "——— |"
Replace-ParaText 3  ". / Beeidigte Ubersetzung aus dem Arabischen"
Replace-ParaText 4  "Republik Tunesien (Wappen der tunesischen Republik)"
Replace-ParaText 5  "Ministerium fiir Hochschulbildung und wissenschaftliche Forschung"
Replace-ParaText 6  "Universitit von Manouba"
Replace-ParaText 7  "Fakuiltit fiir Literaturwissenschaft, Kiinste und Humanwissenschaften"
Replace-ParaText 8  "Das Nationale Zeugnis der Fundamentalen Lizenz (Bachelor)"
Replace-ParaText 9  "¥ Nach Vorlage des Erlasses Nr. 83 des Jahres 1986 vom . September 1986 aber die Errichtung der Fakultat der Literaturwissenschaft von Manouba,"
Replace-ParaText 10 "¥ Nach Vorlage des Gesetzes Nr. 19 des Jahres 2008 vom 25. Februar 2008 aber das Hochschulwesen, insbesondere dessen Artikel 3,"
Replace-ParaText 11 "¥ Nach Vorlage des Erlasses Nr. 1932 des Jahres 2008 vom 02. November 1992 ober die Festiegung der Stelle, die die Unterzeichnung der wissenschaftichen nationalen Hochschulabschlasse zusttindig ist."
Replace-ParaText 12 "“Nach Vorlage des Erfasses Nr. 3123 des Jahres 2008 vom 22. September 2008 uber die Festlegung des allgemeinen Rakmens fbr das Studiensystem und der Bedingungen for den Erwerb des nationalen"
Replace-ParaText 13 "Hochschulabschlusses fur die 1.izenz in den verschiedenen Ausbildungsgebieten, Fachern, Studiengsngen und Fachrichtungen im LMD-System (Lizenz, Master, Doktor},"
Replace-ParaText 14 "¥ Und nach Vorlage der Beratungsprotokolle der Prifungskommissionen des Universitatsjahres 2015-2016,"
Replace-ParaText 15 "wird Frau/ Fraulein: Salma Njema (geboren am 01. 01. 1993 in Monastir, Nationalausweisnummer: 06935513)"
Replace-ParaText 16 "Das Nationale Zeugnis der Fundamentalen Lizenz (Bachelor) in: Fachbereich: Sprachen und Literaturen Hauptfach: Deutsche Sprache, Literatur und Landeskunde mit dem Pradikat: (Ausreichend) erteilt."
Replace-ParaText 17 "Manouba, den 02. 07. 2016"
Replace-ParaText 18 "Der Dekan: Habib Kozdoghli (Unterschrift: Unleserlich) - Dienstsiegel: (Fakultat flir Literaturwissenschaft, Kiinste und"
Replace-ParaText 19 "Humanwissenschaften-In der Mitte: Der Dekan)."
Replace-ParaText 20 "Trockenes Dienstsiegel: (Ministerium fiir Hochschulbildung und wissenschaftliche Forschung- Universitat von Manouba —In der Mitte: Fakultat"

# --- paragraph 21 ("(Signature of the authorised party)", ListBullet) is
#     removed entirely in the new version -----------------------------------
$d.Paragraphs.Item(21).Range.Delete() | Out-Null

# --- remaining two original paragraphs get new text -------------------------
Replace-ParaText 21 "fiir Literaturwissenschaft, Kiinste und Humanwissenschaften von Manouba)"
Replace-ParaText 22 "Hinweis: Das vorliegende Diplom wird nur einmal ausgehindigt."

# --- append the brand new paragraphs at the end of the document ------------
$newParas = @(
    "Auf der Riickseite:",
    "*Stempel des Ministeriums fiir Hochschulbildung und wissenschaftliche Forschung fir die Beglaubigung des Dokumentes: Beglaubigungsvermerk:",
    "Durchsicht erfolgte in der Generaldirektion fiir Hochschulbildung im Ministerium flir Hochschulbildung und wissenschaftliche Forschung. Hiermit bestitigen",
    "wir die Echtheit der Unterschrift des Herm: Der Dekan ohne Verantwortung fiir den Inhalt des vorliegenden Dokuments, Beglaubigungsnummer: 3148, Ort",
    "u. Datum: Tunis, den 22. 01. 2025, Beglaubigungsgebiihr: 5 Dinar, Vizedirektorin der privaten Hochschulbildung: Latifa Ben Abderrahmen Unterschrift (Unleserlich), Siegel des Ministeriums fir Hochschulbildung und wissenschaftliche Forschung (Republik Tunesien - Ministerium fiir",
    "Hochschulbildung und wissenschaftliche Forschung - In der Mitte: Wappen der tunesischen Republik)",
    "*Stempel des AuBenministeriums fiir die Beglaubigung des Dokumentes: Beglaubigungsvermerk: Durchsicht erfolgte im Ministerium fur auswartige",
    "Angelegenheiten. Hiermit bestatigen wir die Echtheit der Unterschrift der Frau: Latifa Ben Abderrahmen, i. A. des Ministers flr Hochschulbildung und",
    "wissenschaftliche Forschung, Ort und Datum: Tunis, den 22. 01. 2025, Beglaubigungsgebihr: 5 Dinar, i. A. des Ministers fur auswirtige Angelegenheiten,",
    "i, A. des Generaldirektors ftir konsularische Angelegenheiten: Hamida Labidi ~ Unterschrift (Unieserlich), Siegel des AuBenministeriums (Republik",
    "Tunesien - Ministerium ftir auswartige Angelegenheiten — In der Mitte: Wappen der tunesischen Republik)",
    "Der Ubersetzung ist eine Kopie des Dokuments angeheftet.",
    "Die Richtigkeit und Vollstindigkeit vorstehender Ubersetzung des mir im Original vorgelegten :",
    "und in arabischer Sprache abgefassten Dokuments wird hiermit bescheinigt. Tunis, den 26. 03. 2025 i",
    "... |"
)

$endRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$endRange.Collapse(0) | Out-Null
foreach ($t in $newParas) {
    $endRange.InsertAfter("`r" + $t)
    $endRange.Collapse(0) | Out-Null
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
